$d = $word.ActiveDocument

function Get-ParaText($para) {
    $t = $para.Range.Text
    # Range.Text includes the trailing paragraph-mark character; strip it.
    if ($t.Length -gt 0 -and $t.Substring($t.Length - 1) -eq [char]13) {
        $t = $t.Substring(0, $t.Length - 1)
    }
    return $t
}

function Set-ParaText($para, $text) {
    $r = $para.Range
    $r.MoveEnd(1, -1) | Out-Null   # exclude trailing paragraph mark
    $r.Text = $text
}

# Find the given literal text anywhere in the document and overwrite it by
# direct Range.Text assignment (NOT the Find.Execute "Replace" argument,
# which runs the inserted string through AutoCorrect/smart-quotes).
function Replace-ExactText($doc, $oldText, $newText) {
    $rng = $doc.Content
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $oldText"
    }
    $rng.Text = $newText
}

# --- Capture the "before" text of every whole paragraph involved in the
# --- content rotation, before any mutation happens (index-based reads,
# --- unaffected by ordering of the Find calls below).

$oldP6  = Get-ParaText $d.Paragraphs.Item(6)   # Objetivos body
$oldP8  = Get-ParaText $d.Paragraphs.Item(8)   # Docente(s) body
$oldP10 = Get-ParaText $d.Paragraphs.Item(10)  # Programa resumido body
$oldP12 = Get-ParaText $d.Paragraphs.Item(12)  # Programa body
$oldP16 = Get-ParaText $d.Paragraphs.Item(16)  # Bibliografia body

# Paragraph 14 (Avaliacao body) holds three labelled segments in separate
# runs: "Metodo: ", "Criterio: ", "Norma de recuperacao: ", each followed
# by a value run. These values are unique, known strings in the original
# document - use them directly as Find anchors (preserves the bold label
# runs untouched).
$oldMetodo = "Aplicação de prova(s) e relatório(s)."
$oldCriterio = "A média do período será definida pelo professor da disciplina. Alunos com média final igual ou superior a 5,0 estarão aprovados, desde que tenham freqüência mínima de 70% (regimental). Alunos com média inferior a 3,0 e/ou freqüência inferior a 70% estarão reprovados (regimental). Alunos com média superior ou igual a 3,0 e inferior a 5,0 e que tenham freqüência mínima de 70% serão submetidos ao período de recuperação (regimental)."
$oldNorma = "A média final após a recuperação para a disciplina será a média aritmética entre a média do período e a nota da recuperação. Durante o período de recuperação, poderá ser marcada uma aula com a finalidade de sanar dúvidas e/ou revisar conceitos fundamentais. Em data posterior os alunos serão submetidos a uma prova de recuperação"

$NL = [char]11   # manual line break (<w:br/>) as used inside Range.Text

$bibText = "1) FOUST, Alan S.; WENZEL, Leonard A.; CLUMP, Curtis W.; MAUS, Louis; ANDERSEN, L. Bryce. Princípios das operações unitárias. Rio de Janeiro: Guanabara Dois/LTC, 1982." + $NL + `
"2) GEANKOPLIS, Christie John. Transport Processes and Separation Process Principles. New York: Prentice Hall, 2003." + $NL + `
"3) COUPER, James R.; PENNEY, W. Roy; FAIR, James R.; WALAS, Stanley M. Chemical Process Equipment: Selection and Design. Amsterdam: Elsevier, 2005." + $NL + `
"4) FOGLER, H. S. Elementos de engenharia das reações químicas. 3.ed. Rio de Janeiro: LTC Editora, 2002." + $NL + `
"5) LEVENSPIEL, O. Chemical Reaction Engineering. 3rd.ed. New York: John Wiley & Sons, 1998." + $NL + `
"6) PERRY, Robert H.; GREEN, Don W. Perry's Chemical Engineers' Handbook. 8th.ed. New York: McGraw-Hill, 2008."

# --- Apply the rotation: each destination receives the text that used to
# --- live one step earlier in the cycle:
# ---   P6 <- P10 <- P12 <- Metodo <- Criterio <- Norma <- P16 <- P8 <- P6
#
# The three Find/replace calls below must run first, while oldMetodo /
# oldCriterio / oldNorma are each still unique in the document (before
# oldMetodo gets written into paragraph 12 further down).

Replace-ExactText $d $oldNorma $bibText
Replace-ExactText $d $oldCriterio $oldNorma
Replace-ExactText $d $oldMetodo $oldCriterio

Set-ParaText $d.Paragraphs.Item(12) $oldMetodo
Set-ParaText $d.Paragraphs.Item(10) $oldP12
Set-ParaText $d.Paragraphs.Item(6) $oldP10
Set-ParaText $d.Paragraphs.Item(16) $oldP8
Set-ParaText $d.Paragraphs.Item(8) $oldP6
